# Apply the "cousybo01" edit:
#  - shared-string header rename "Team" -> "Tm" (affects every sheet that
#    uses that header, i.e. per_game, per_minute, advanced)
#  - on the "advanced" sheet: insert a new blank column before the old
#    "RSPS" column (U), give the new column a single non-breaking-space
#    header, and tweak a handful of PER (column I) values
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "Team" header to "Tm" everywhere it is used.
# ---------------------------------------------------------------------
foreach ($sheetName in @("per_game", "per_minute", "advanced")) {
    $sheet = $wb.Sheets.Item($sheetName)
    if ($sheet.Range("D1").Value2 -eq "Team") {
        $sheet.Range("D1").Value = "Tm"
    }
}

# ---------------------------------------------------------------------
# 2. "advanced" sheet: insert a new column U (old U -> V), and give the
#    new column a single non-breaking-space header in row 1.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("advanced")
$ws.Columns("U:U").Insert()
$ws.Range("U1").Value = [char]0x00A0

# ---------------------------------------------------------------------
# 3. Update the handful of recalculated PER (column I) values.
# ---------------------------------------------------------------------
$ws.Range("I5").Value = 21.5
$ws.Range("I7").Value = 3.899999999999999
$ws.Range("I8").Value = 21.7
$ws.Range("I10").Value = 1
$ws.Range("I14").Value = 21.7
$ws.Range("I16").Value = -2
$ws.Range("I17").Value = 20.3
$ws.Range("I19").Value = 11.3
$ws.Range("I20").Value = 21.2
$ws.Range("I22").Value = -4.099999999999998
$ws.Range("I35").Value = 18.5
$ws.Range("I37").Value = -3
$ws.Range("I43").Value = 19.9
$ws.Range("I45").Value = -2.5
$ws.Range("I48").Value = 19.9
$ws.Range("I50").Value = -2.5
